$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 61
$ws.Range("A60").Copy()
$ws.Range("A61").PasteSpecial(-4122)
$ws.Range("A61").Value = 45629
$ws.Range("B61").Value = 361
$ws.Range("D61").Value = 2
$ws.Range("E61").Value = "some miss entry so observer the trade on Sunday"

# Row 63 text written before row 62 text so the shared-string table gains
# entries in the same order as the source workbook (177, 178, 179).
$ws.Range("E63").Value = "I knew there was strong selling pressure as per oi but I went against the data so I boooked loss in 1 trade"

# Row 62
$ws.Range("A60").Copy()
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("A62").Value = 45630
$ws.Range("B62").Value = 393
$ws.Range("D62").Value = 1
$ws.Range("E62").Value = "it was a good trade indeed"

# Row 63 (remaining cells)
$ws.Range("A60").Copy()
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("A63").Value = 45631
$ws.Range("B63").Value = 512
$ws.Range("D63").Value = 2

$excel.CutCopyMode = $false

$ws.Range("C63").Select()
